
$d = $word.ActiveDocument

# Find the placeholder paragraph (Day 7 bullet currently containing just the
# ellipsis character "…") that needs to be filled in.
$ellipsis = [string][char]0x2026
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $ellipsis) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the placeholder '…' paragraph"
}

# Replace that single placeholder paragraph with two fully-formed list items:
#   1) "Formulated logical programming steps for the chosen tasks" (with the
#      trailing word wrapped in gramStart/gramEnd proofErr markers, matching
#      the grammar-check artifacts Word leaves on similar bullets elsewhere
#      in this document).
#   2) "Started implementation of both tasks parallelly"
# Both keep the same ListParagraph / numbered-list (numId 1, ilvl 0) format
# as the paragraph being replaced.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Formulated logical programming steps for the chosen </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>tasks</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Started implementation of both tasks parallelly</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml)

# Re-locate the just-inserted "Started implementation..." paragraph and add a
# single, unformatted blank paragraph right after it (ahead of the block of
# existing blank paragraphs that already trails this section). Inserting it
# as the paragraph *before* the next (already-blank, style-less) paragraph
# lets it pick up that plain formatting instead of the ListParagraph/bullet
# formatting it would inherit from InsertParagraphAfter on the bullet above.
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Started implementation of both tasks parallelly*") {
        $following = $para.Next()
        $following.Range.InsertParagraphBefore()
        break
    }
}
